$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")

# ALC!row 4
$ws_ALC.Range("H4").Value = 350
$ws_ALC.Range("I4").Value = 600
$ws_ALC.Range("K4").Value = 600
$ws_ALC.Range("M4").Value = -486
# ALC!row 6
$ws_ALC.Range("H6").Value = 1044
$ws_ALC.Range("I6").Value = 66
$ws_ALC.Range("J6").Value = 3000
$ws_ALC.Range("K6").Value = 198
$ws_ALC.Range("L6").Value = 9000
$ws_ALC.Range("M6").Value = -86
$ws_ALC.Range("N6").Value = -9224
# ALC!row 8
$ws_ALC.Range("H8").Value = 10.166667
$ws_ALC.Range("I8").Value = 6.2
$ws_ALC.Range("J8").Value = 30
$ws_ALC.Range("K8").Value = 18.6
$ws_ALC.Range("L8").Value = 90
$ws_ALC.Range("M8").Value = 120.4
$ws_ALC.Range("N8").Value = -368
# ALC!row 10
$ws_ALC.Range("H10").Value = 0
$ws_ALC.Range("J10").Value = 0
$ws_ALC.Range("L10").Value = 0
$ws_ALC.Range("N10").ClearContents()
# ALC!row 12
$ws_ALC.Range("H12").Value = 465.55554
$ws_ALC.Range("I12").Value = 395
$ws_ALC.Range("J12").Value = 485.7143
$ws_ALC.Range("K12").Value = 395
$ws_ALC.Range("L12").Value = 485.7143
$ws_ALC.Range("M12").Value = -225
$ws_ALC.Range("N12").Value = -825.7143
# ALC!row 18
$ws_ALC.Range("H18").Value = 2000
$ws_ALC.Range("I18").Value = 2000
$ws_ALC.Range("K18").Value = 2000
$ws_ALC.Range("M18").Value = -1716
# ALC!row 31
$ws_ALC.Range("H31").Value = 3036
$ws_ALC.Range("J31").Value = 15000
$ws_ALC.Range("L31").Value = 45000
$ws_ALC.Range("N31").Value = -45460
# ALC!row 38
$ws_ALC.Range("H38").Value = 1915.8125
$ws_ALC.Range("I38").Value = 408
$ws_ALC.Range("K38").Value = 1224
$ws_ALC.Range("M38").Value = -852
# ALC!row 88
$ws_ALC.Range("H88").Value = 1662.5834
$ws_ALC.Range("I88").Value = 441.66666
$ws_ALC.Range("J88").Value = 2069.5557
$ws_ALC.Range("K88").Value = 441.66666
$ws_ALC.Range("L88").Value = 2069.5557
$ws_ALC.Range("M88").Value = -35.66665999999998
$ws_ALC.Range("N88").Value = -2881.5557
# ALC!row 91
$ws_ALC.Range("H91").Value = 1662.5834
$ws_ALC.Range("I91").Value = 441.66666
$ws_ALC.Range("J91").Value = 2069.5557
$ws_ALC.Range("K91").Value = 441.66666
$ws_ALC.Range("L91").Value = 2069.5557
$ws_ALC.Range("M91").Value = 962.33334
$ws_ALC.Range("N91").Value = -4877.5557
# ALC!row 107
$ws_ALC.Range("H107").Value = 157.42857
$ws_ALC.Range("I107").Value = 157.42857
$ws_ALC.Range("K107").Value = 157.42857
$ws_ALC.Range("M107").Value = 1762.57143
# ALC!row 137
$ws_ALC.Range("H137").Value = 1682.2858
$ws_ALC.Range("I137").Value = 1225
$ws_ALC.Range("K137").Value = 3675
$ws_ALC.Range("M137").Value = -1125
# ARM!row 5
$ws_ARM.Range("H5").Value = 490
$ws_ARM.Range("I5").Value = 487.5
$ws_ARM.Range("J5").Value = 500
$ws_ARM.Range("K5").Value = 487.5
$ws_ARM.Range("L5").Value = 500
$ws_ARM.Range("M5").Value = -375.5
$ws_ARM.Range("N5").Value = -724
# ARM!row 26
$ws_ARM.Range("H26").Value = 1567.8334
$ws_ARM.Range("I26").Value = 1567.8334
$ws_ARM.Range("K26").Value = 1567.8334
$ws_ARM.Range("M26").Value = -1237.8334
# ARM!row 32
$ws_ARM.Range("H32").Value = 9904.308000000001
$ws_ARM.Range("I32").Value = 8250.546
$ws_ARM.Range("K32").Value = 8250.546
$ws_ARM.Range("M32").Value = -7963.546
# ARM!row 74
$ws_ARM.Range("H74").Value = 6497.393
$ws_ARM.Range("I74").Value = 6596.6
$ws_ARM.Range("J74").Value = 6249.375
$ws_ARM.Range("K74").Value = 6596.6
$ws_ARM.Range("L74").Value = 6249.375
$ws_ARM.Range("M74").Value = -5722.6
$ws_ARM.Range("N74").Value = -7997.375
# ARM!row 77
$ws_ARM.Range("H77").Value = 6497.393
$ws_ARM.Range("I77").Value = 6596.6
$ws_ARM.Range("J77").Value = 6249.375
$ws_ARM.Range("K77").Value = 32983
$ws_ARM.Range("L77").Value = 31246.875
$ws_ARM.Range("M77").Value = -28615
$ws_ARM.Range("N77").Value = -39982.875
# ARM!row 132
$ws_ARM.Range("H132").Value = 2447.8333
$ws_ARM.Range("J132").Value = 1800
$ws_ARM.Range("L132").Value = 5400
$ws_ARM.Range("N132").Value = -10460
# BSM!row 4
$ws_BSM.Range("H4").Value = 490
$ws_BSM.Range("I4").Value = 487.5
$ws_BSM.Range("J4").Value = 500
$ws_BSM.Range("K4").Value = 487.5
$ws_BSM.Range("L4").Value = 500
$ws_BSM.Range("M4").Value = -372.5
$ws_BSM.Range("N4").Value = -730
# BSM!row 22
$ws_BSM.Range("H22").Value = 73
$ws_BSM.Range("I22").Value = 73
$ws_BSM.Range("K22").Value = 73
$ws_BSM.Range("M22").Value = 100
# BSM!row 80
$ws_BSM.Range("H80").Value = 131.47058
$ws_BSM.Range("I80").Value = 76.166664
$ws_BSM.Range("J80").Value = 161.63637
$ws_BSM.Range("K80").Value = 76.166664
$ws_BSM.Range("L80").Value = 161.63637
$ws_BSM.Range("M80").Value = 921.833336
$ws_BSM.Range("N80").Value = -2157.63637
# BSM!row 83
$ws_BSM.Range("H83").Value = 131.47058
$ws_BSM.Range("I83").Value = 76.166664
$ws_BSM.Range("J83").Value = 161.63637
$ws_BSM.Range("K83").Value = 380.83332
$ws_BSM.Range("L83").Value = 808.1818499999999
$ws_BSM.Range("M83").Value = 4611.16668
$ws_BSM.Range("N83").Value = -10792.18185
# BSM!row 107
$ws_BSM.Range("H107").Value = 6444.9165
$ws_BSM.Range("J107").Value = 8875
$ws_BSM.Range("L107").Value = 8875
$ws_BSM.Range("N107").Value = -12715
# CRP!row 22
$ws_CRP.Range("H22").Value = 2025.8
$ws_CRP.Range("I22").Value = 1120.7778
$ws_CRP.Range("K22").Value = 1120.7778
$ws_CRP.Range("M22").Value = -770.7778000000001
# CUL!row 17
$ws_CUL.Range("H17").Value = 71.5
$ws_CUL.Range("I17").Value = 25.8
$ws_CUL.Range("J17").Value = 300
$ws_CUL.Range("K17").Value = 77.40000000000001
$ws_CUL.Range("L17").Value = 900
$ws_CUL.Range("M17").Value = 91.59999999999999
$ws_CUL.Range("N17").Value = -1238
# CUL!row 40
$ws_CUL.Range("H40").Value = 28.75
$ws_CUL.Range("J40").Value = 0
$ws_CUL.Range("L40").Value = 0
$ws_CUL.Range("N40").ClearContents()
# GSM!row 2
$ws_GSM.Range("H2").Value = 40.090908
$ws_GSM.Range("I2").Value = 32
$ws_GSM.Range("J2").Value = 43.125
$ws_GSM.Range("K2").Value = 32
$ws_GSM.Range("L2").Value = 43.125
$ws_GSM.Range("M2").Value = 81
$ws_GSM.Range("N2").Value = -269.125
# GSM!row 17
$ws_GSM.Range("H17").Value = 0
$ws_GSM.Range("J17").Value = 0
$ws_GSM.Range("L17").Value = 0
$ws_GSM.Range("N17").ClearContents()
# GSM!row 21
$ws_GSM.Range("H21").Value = 0
$ws_GSM.Range("J21").Value = 0
$ws_GSM.Range("L21").Value = 0
$ws_GSM.Range("N21").ClearContents()
# GSM!row 30
$ws_GSM.Range("H30").Value = 0
$ws_GSM.Range("J30").Value = 0
$ws_GSM.Range("L30").Value = 0
$ws_GSM.Range("N30").ClearContents()
# GSM!row 62
$ws_GSM.Range("H62").Value = 0
$ws_GSM.Range("I62").Value = 0
$ws_GSM.Range("K62").Value = 0
$ws_GSM.Range("M62").ClearContents()
# GSM!row 65
$ws_GSM.Range("H65").Value = 0
$ws_GSM.Range("I65").Value = 0
$ws_GSM.Range("K65").Value = 0
$ws_GSM.Range("M65").ClearContents()
# GSM!row 70
$ws_GSM.Range("H70").Value = 3328.4167
$ws_GSM.Range("I70").Value = 3394.2
$ws_GSM.Range("J70").Value = 2999.5
$ws_GSM.Range("K70").Value = 3394.2
$ws_GSM.Range("L70").Value = 2999.5
$ws_GSM.Range("M70").Value = -3124.2
$ws_GSM.Range("N70").Value = -3539.5
# GSM!row 73
$ws_GSM.Range("H73").Value = 3328.4167
$ws_GSM.Range("I73").Value = 3394.2
$ws_GSM.Range("J73").Value = 2999.5
$ws_GSM.Range("K73").Value = 3394.2
$ws_GSM.Range("L73").Value = 2999.5
$ws_GSM.Range("M73").Value = -2458.2
$ws_GSM.Range("N73").Value = -4871.5
# GSM!row 107
$ws_GSM.Range("H107").Value = 362
$ws_GSM.Range("J107").Value = 463.33334
$ws_GSM.Range("L107").Value = 463.33334
$ws_GSM.Range("N107").Value = -4303.33334
# LTW!row 35
$ws_LTW.Range("H35").Value = 2580.7144
$ws_LTW.Range("I35").Value = 1698.3334
$ws_LTW.Range("J35").Value = 3242.5
$ws_LTW.Range("K35").Value = 1698.3334
$ws_LTW.Range("L35").Value = 3242.5
$ws_LTW.Range("M35").Value = -1362.3334
$ws_LTW.Range("N35").Value = -3914.5
# LTW!row 46
$ws_LTW.Range("H46").Value = 4250.4375
$ws_LTW.Range("I46").Value = 3535.3333
$ws_LTW.Range("J46").Value = 4679.5
$ws_LTW.Range("K46").Value = 3535.3333
$ws_LTW.Range("L46").Value = 4679.5
$ws_LTW.Range("M46").Value = -3347.3333
$ws_LTW.Range("N46").Value = -5055.5
# LTW!row 133
$ws_LTW.Range("H133").Value = 0
$ws_LTW.Range("J133").Value = 0
$ws_LTW.Range("L133").Value = 0
$ws_LTW.Range("N133").ClearContents()
